$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# 1) First three rows' values ("99.99", "0", "69") are replaced with "0M"
$t.Cell(1,1).Range.Text = "0M"
$t.Cell(2,1).Range.Text = "0M"
$t.Cell(3,1).Range.Text = "0M"

# 2) Insert 10 new single-value rows right after row 3 (i.e. before the
#    row that currently holds the old row 4's "0" value). Insert them in
#    reverse order, always before the same reference row, so they land in
#    the desired forward order.
$values = @("104", "0.00003", "0.00011", "0.00006", "0.00002", "0.00007", "0.00007", "0.00011", "0.00463", "100.0")
$refRow = $t.Rows.Item(4)
for ($i = $values.Length - 1; $i -ge 0; $i--) {
    $newRow = $t.Rows.Add($refRow)
    $newRow.Cells(1).Range.Text = $values[$i]
    $refRow = $newRow
}

# 3) The three tail rows that used to hold tab-separated value lists are
#    collapsed down to a single value each, re-using the values that used
#    to live in the original rows 1-3.
$t.Cell(44,1).Range.Text = "99.99"
$t.Cell(45,1).Range.Text = "0"
$t.Cell(46,1).Range.Text = "69"
